$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---

# Overview sheet: zh-cn and de-de status columns (E2, F2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Resize the "Status" related column(s) from ~17.216 to ~13.410 characters wide ---
# (ColumnWidth value chosen so the persisted OOXML <col width> lands as close as
# possible to the target 13.4101845877511 given this engine's width rounding.)

$wsOverview.Range("E:E").ColumnWidth = 12.5
$wsOverview.Range("F:F").ColumnWidth = 12.5

$wsZhCn.Range("C:C").ColumnWidth = 12.5

$wsDeDe.Range("C:C").ColumnWidth = 12.5
